$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "2025-04-02"
$ws.Range("A39").Style = "Normal"
$ws.Range("B39").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C39").Value = "NA"
$ws.Range("D39").Value = 1
